$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.479329347610474
$ws.Range("B1").Value = 2.011158466339111
$ws.Range("C1").Value = 2.358525514602661
$ws.Range("D1").Value = 2.815152645111084
$ws.Range("E1").Value = 2.72649621963501
